$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "41.067.58"
$ws.Range("E2").Value = "  -2.40%  "

# Row 3
$ws.Range("D3").Value = "2.154.11"
$ws.Range("E3").Value = "  -2.94%  "

# Row 4
$ws.Range("E4").Value = "  +0.17%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "235.43"
$ws.Range("E5").Value = "  -2.57%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.600"
$ws.Range("E6").Value = "  -4.72%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "68.55"
$ws.Range("E7").Value = "  -6.56%  "

# Row 8
$ws.Range("E8").Value = "  +0.06%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.560"
$ws.Range("E9").Value = "  -7.94%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "38.08"
$ws.Range("E10").Value = "  -10.51%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0897"
$ws.Range("E11").Value = "  -6.11%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "54.13"
$ws.Range("E12").Value = "  -5.79%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0995"
$ws.Range("E13").Value = "  -4.05%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.60"
$ws.Range("E14").Value = "  -6.93%  "

# Row 15
$ws.Range("D15").Value = "2.479.74"
$ws.Range("E15").Value = "  -2.84%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.26"
$ws.Range("E16").Value = "  -0.41%  "

# Row 17
$ws.Range("D17").Value = "2.131.49"
$ws.Range("E17").Value = "  -3.99%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.776"
$ws.Range("E18").Value = "  -7.48%  "

# Row 19
$ws.Range("D19").Value = "40.906.85"
$ws.Range("E19").Value = "  -2.48%  "

# Row 20
$ws.Range("D20").Value = "0.0₃0983"
$ws.Range("E20").Value = "  -8.48%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "69.44"
$ws.Range("E21").Value = "  -4.84%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.71"
$ws.Range("E22").Value = "  -8.36%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "222.89"
$ws.Range("E23").Value = "  -3.09%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.34"
$ws.Range("E24").Value = "  -13.85%  "

# Row 25
$ws.Range("E25").Value = "  +0.06%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.85"
$ws.Range("E26").Value = "  -11.90%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.52"
$ws.Range("E27").Value = "  -11.15%  "

# Row 28
$ws.Range("B28").Value = "PancakeSwap"
$ws.Range("C28").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.17"
$ws.Range("E28").Value = "  -4.55%  "

# Row 29
$ws.Range("B29").Value = "WEMIXToken"
$ws.Range("C29").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "3.19"
$ws.Range("E29").Value = "  -11.58%  "

# Row 30
$ws.Range("E30").Value = "  -1.62%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "168.09"
$ws.Range("E31").Value = "  +0.31%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "19.57"
$ws.Range("E32").Value = "  -4.63%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "29.82"
$ws.Range("E33").Value = "  +0.56%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0745"
$ws.Range("E34").Value = "  -6.76%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.00"
$ws.Range("E35").Value = "  -9.96%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.119"
$ws.Range("E36").Value = "  -4.88%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.101"
$ws.Range("E37").Value = "  -8.24%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.04"
$ws.Range("E38").Value = "  -5.64%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0275"
$ws.Range("E39").Value = "  -8.69%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.03"
$ws.Range("E40").Value = "  -4.50%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "11.30"
$ws.Range("E41").Value = "  -16.90%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.26"
$ws.Range("E42").Value = "  -6.80%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "56.82"
$ws.Range("E43").Value = "  -13.95%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.183"
$ws.Range("E44").Value = "  -7.80%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.14"
$ws.Range("E45").Value = "  -7.25%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0949"
$ws.Range("E46").Value = "  -5.49%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "95.74"
$ws.Range("E47").Value = "  -9.04%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.06"
$ws.Range("E48").Value = "  -4.95%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.10"
$ws.Range("E49").Value = "  -5.47%  "

# Row 50
$ws.Range("B50").Value = "HuobiToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.60"
$ws.Range("E50").Value = "  -3.24%  "

# Row 51
$ws.Range("B51").Value = "NEARProtocol"
$ws.Range("C51").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.13"
$ws.Range("E51").Value = "  -11.94%  "
